$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 47, shifting the existing rows 47:60 down to 49:62.
$ws.Rows("47:48").Insert(-4121)

# --- Row 47 (new weekly entry) ---
$ws.Range("A47").Value = 7
$ws.Range("B47").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C47").Value = "Ñuble"
$ws.Range("D47").Value = 45124
$ws.Range("E47").Value = 16
$ws.Range("F47").Value = "Fruta"
$ws.Range("G47").Value = 100104
$ws.Range("H47").Value = "Frutos de pepita"
$ws.Range("I47").Value = 100104003
$ws.Range("J47").Value = "Membrillo"
$ws.Range("K47").Value = "Champion"
$ws.Range("L47").Value = "Primera"
$ws.Range("M47").Value = 80
$ws.Range("N47").Value = 10000
$ws.Range("O47").Value = 10000
$ws.Range("P47").Value = 10000
$ws.Range("Q47").Value = "`$/bandeja 18 kilos granel"
$ws.Range("R47").Value = "Región de O'Higgins"
$ws.Range("S47").Value = 556
$ws.Range("T47").Value = 18

# --- Row 48 (new weekly entry) ---
$ws.Range("A48").Value = 7
$ws.Range("B48").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C48").Value = "Ñuble"
$ws.Range("D48").Value = 45124
$ws.Range("E48").Value = 16
$ws.Range("F48").Value = "Fruta"
$ws.Range("G48").Value = 100104
$ws.Range("H48").Value = "Frutos de pepita"
$ws.Range("I48").Value = 100104003
$ws.Range("J48").Value = "Membrillo"
$ws.Range("K48").Value = "Champion"
$ws.Range("L48").Value = "Segunda"
$ws.Range("M48").Value = 80
$ws.Range("N48").Value = 8000
$ws.Range("O48").Value = 8000
$ws.Range("P48").Value = 8000
$ws.Range("Q48").Value = "`$/bandeja 18 kilos granel"
$ws.Range("R48").Value = "Región de O'Higgins"
$ws.Range("S48").Value = 444
$ws.Range("T48").Value = 18
